# Update Excel workbook for conversion (matches commit "Updating excel file for conversion")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 4-11 (X/Y coordinate values recalculated) ---
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 100

$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 100

$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 125

$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 125

$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 75

# Rows 10-11 change descriptor from WALL to the new WINDOW string
$ws.Range("B10").Value = "WINDOW"
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 75

$ws.Range("B11").Value = "WINDOW"
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 25

# --- Append new rows 12-14 (Exterior / WALL entries) ---
$ws.Range("A12").Value = "Exterior"
$ws.Range("B12").Value = "WALL"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 25

$ws.Range("A13").Value = "Exterior"
$ws.Range("B13").Value = "WALL"
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 0

$ws.Range("A14").Value = "Exterior"
$ws.Range("B14").Value = "WALL"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0

# --- Update the view: selected cell moves to the new last row (F14) ---
$ws.Range("F14").Select()
